$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 25802.625
$ws.Range("J17").Value = 27862.297
$ws.Range("L17").Value = 83586.891
$ws.Range("N17").Value = -83922.891

$ws.Range("H80").Value = 5173770
$ws.Range("I80").Value = 3270071.8
$ws.Range("J80").Value = 6805511.5
$ws.Range("K80").Value = 9810215.399999999
$ws.Range("L80").Value = 20416534.5
$ws.Range("M80").Value = -9809217.399999999
$ws.Range("N80").Value = -20418530.5

$ws.Range("H83").Value = 5173770
$ws.Range("I83").Value = 3270071.8
$ws.Range("J83").Value = 6805511.5
$ws.Range("K83").Value = 29430646.2
$ws.Range("L83").Value = 61249603.5
$ws.Range("M83").Value = -29425654.2
$ws.Range("N83").Value = -61259587.5

$ws.Range("H86").Value = 4175.2607
$ws.Range("I86").Value = 2132.2307
$ws.Range("J86").Value = 6831.2
$ws.Range("K86").Value = 2132.2307
$ws.Range("L86").Value = 6831.2
$ws.Range("M86").Value = -1009.2307
$ws.Range("N86").Value = -9077.200000000001

$ws.Range("H89").Value = 4175.2607
$ws.Range("I89").Value = 2132.2307
$ws.Range("J89").Value = 6831.2
$ws.Range("K89").Value = 10661.1535
$ws.Range("L89").Value = 34156
$ws.Range("M89").Value = -5045.1535
$ws.Range("N89").Value = -45388

$ws.Range("H106").Value = 9279.6875
$ws.Range("I106").Value = 8728.846
$ws.Range("K106").Value = 8728.846
$ws.Range("M106").Value = -8097.846

$ws.Range("H137").Value = 928697.2
$ws.Range("I137").Value = 1305.6
$ws.Range("J137").Value = 1443914.8
$ws.Range("K137").Value = 3916.8
$ws.Range("L137").Value = 4331744.4
$ws.Range("M137").Value = -1366.8
$ws.Range("N137").Value = -4336844.4

$ws.Range("H140").Value = 108999
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4299.625
$ws.Range("I32").Value = 4219.68
$ws.Range("K32").Value = 4219.68
$ws.Range("M32").Value = -3932.68

$ws.Range("H61").Value = 16159445
$ws.Range("I61").Value = 22227708
$ws.Range("K61").Value = 22227708
$ws.Range("M61").Value = -22227496

$ws.Range("H74").Value = 1089001.8
$ws.Range("I74").Value = 1251711.5
$ws.Range("J74").Value = 4270.3335
$ws.Range("K74").Value = 1251711.5
$ws.Range("L74").Value = 4270.3335
$ws.Range("M74").Value = -1250837.5
$ws.Range("N74").Value = -6018.3335

$ws.Range("H77").Value = 1089001.8
$ws.Range("I77").Value = 1251711.5
$ws.Range("J77").Value = 4270.3335
$ws.Range("K77").Value = 6258557.5
$ws.Range("L77").Value = 21351.6675
$ws.Range("M77").Value = -6254189.5
$ws.Range("N77").Value = -30087.6675

$ws.Range("H132").Value = 2504378.2
$ws.Range("I132").Value = 4160.8857
$ws.Range("K132").Value = 12482.6571
$ws.Range("M132").Value = -9952.6571

$ws.Range("H136").Value = 16159445
$ws.Range("I136").Value = 22227708
$ws.Range("K136").Value = 66683124
$ws.Range("M136").Value = -66680574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3007.5
$ws.Range("J99").Value = 4799.6665
$ws.Range("L99").Value = 4799.6665
$ws.Range("N99").Value = -7795.6665

$ws.Range("H134").Value = 3126750.5
$ws.Range("I134").Value = 1566.8889
$ws.Range("K134").Value = 4700.6667
$ws.Range("M134").Value = -2165.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29711448
$ws.Range("I31").Value = 50002708
$ws.Range("J31").Value = 723936.8
$ws.Range("K31").Value = 50002708
$ws.Range("L31").Value = 723936.8
$ws.Range("M31").Value = -50002413
$ws.Range("N31").Value = -724526.8

$ws.Range("H34").Value = 29711448
$ws.Range("I34").Value = 50002708
$ws.Range("J34").Value = 723936.8
$ws.Range("K34").Value = 50002708
$ws.Range("L34").Value = 723936.8
$ws.Range("M34").Value = -50002506
$ws.Range("N34").Value = -724340.8

$ws.Range("H122").Value = 4419.5713
$ws.Range("J122").Value = 4066.5715
$ws.Range("L122").Value = 12199.7145
$ws.Range("N122").Value = -17099.7145

$ws.Range("H132").Value = 2345.7334
$ws.Range("I132").Value = 2138.12
$ws.Range("K132").Value = 6414.36
$ws.Range("M132").Value = -3884.36

$ws.Range("H134").Value = 2758.4167
$ws.Range("I134").Value = 2524.2778
$ws.Range("J134").Value = 3460.8333
$ws.Range("K134").Value = 7572.8334
$ws.Range("L134").Value = 10382.4999
$ws.Range("M134").Value = -5037.8334
$ws.Range("N134").Value = -15452.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1598
$ws.Range("I5").Value = 760
$ws.Range("J5").Value = 3832.6667
$ws.Range("K5").Value = 2280
$ws.Range("L5").Value = 11498.0001
$ws.Range("M5").Value = -2168
$ws.Range("N5").Value = -11722.0001

$ws.Range("H51").Value = 12945.5
$ws.Range("I51").Value = 1724.5
$ws.Range("J51").Value = 24166.5
$ws.Range("K51").Value = 5173.5
$ws.Range("L51").Value = 72499.5
$ws.Range("M51").Value = -4713.5
$ws.Range("N51").Value = -73419.5

$ws.Range("H68").Value = 4982.6665
$ws.Range("J68").Value = 4982.6665
$ws.Range("L68").Value = 14947.9995
$ws.Range("N68").Value = -16569.9995

$ws.Range("H71").Value = 4982.6665
$ws.Range("J71").Value = 4982.6665
$ws.Range("L71").Value = 44843.9985
$ws.Range("N71").Value = -52955.9985

$ws.Range("H80").Value = 41672416

$ws.Range("H83").Value = 41672416

$ws.Range("H115").Value = 17029
$ws.Range("I115").Value = 725
$ws.Range("J115").Value = 33333
$ws.Range("K115").Value = 2175
$ws.Range("L115").Value = 99999
$ws.Range("M115").Value = -1000
$ws.Range("N115").Value = -102349

$ws.Range("H135").Value = 1598
$ws.Range("I135").Value = 760
$ws.Range("J135").Value = 3832.6667
$ws.Range("K135").Value = 6840
$ws.Range("L135").Value = 34494.0003
$ws.Range("M135").Value = -4305
$ws.Range("N135").Value = -39564.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5656
$ws.Range("I43").Value = 5656
$ws.Range("K43").Value = 5656
$ws.Range("M43").Value = -5505

$ws.Range("H57").Value = 46666.668
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 46666.668
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 46666.668
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -48306.668

$ws.Range("H70").Value = 10403.714
$ws.Range("I70").Value = 8867.833000000001
$ws.Range("K70").Value = 8867.833000000001
$ws.Range("M70").Value = -8597.833000000001

$ws.Range("H73").Value = 10403.714
$ws.Range("I73").Value = 8867.833000000001
$ws.Range("K73").Value = 8867.833000000001
$ws.Range("M73").Value = -7931.833000000001

$ws.Range("H97").Value = 877.1852
$ws.Range("I97").Value = 838.65216
$ws.Range("K97").Value = 838.65216
$ws.Range("M97").Value = -342.65216

$ws.Range("H132").Value = 11934956
$ws.Range("I132").Value = 3500.6
$ws.Range("K132").Value = 10501.8
$ws.Range("M132").Value = -7971.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3060.1448
$ws.Range("I122").Value = 2993.739
$ws.Range("J122").Value = 3714.7144
$ws.Range("K122").Value = 8981.217000000001
$ws.Range("L122").Value = 11144.1432
$ws.Range("M122").Value = -6531.217000000001
$ws.Range("N122").Value = -16044.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H100").Value = 1666.3334
$ws.Range("I100").Value = 1666.3334
$ws.Range("K100").Value = 3332.6668
$ws.Range("M100").Value = -2791.6668

$ws.Range("H107").Value = 3011.5625
$ws.Range("I107").Value = 1303.3043
$ws.Range("K107").Value = 3909.9129
$ws.Range("M107").Value = -1989.9129

$ws.Range("H113").Value = 999.3333
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340

$ws.Range("H132").Value = 346241.97
$ws.Range("I132").Value = 1435.7826
$ws.Range("K132").Value = 4307.3478
$ws.Range("M132").Value = -1777.3478

$ws.Range("H136").Value = 272492.1
$ws.Range("I136").Value = 2291.6287
$ws.Range("K136").Value = 6874.886100000001
$ws.Range("M136").Value = -4324.886100000001
